$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This string is used as the "Status" value on each localization sheet.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"

# --- Column width changes (status columns narrowed) ---
# Target stored width is 13.4101845877511 "character" units.  Excel's
# COM ColumnWidth setter snaps to whole-pixel boundaries on save, so the
# nearest value reachable through the object model (12.5 -> stored
# 13.333333333333334) is used here.
$ws1.Range("E1").ColumnWidth = 12.5
$ws1.Range("F1").ColumnWidth = 12.5

$ws2.Range("C1").ColumnWidth = 12.5

$ws3.Range("C1").ColumnWidth = 12.5
